$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for columns B, C, D, F across rows 2-10 (row 1 is header, unchanged)
$data = @(
    @{ B = "NSE:ABBOTINDIA"; C = "NSE:ACI";        D = "NSE:ANGELONE"; F = "NSE:M&MFIN" },
    @{ B = "NSE:CANTABIL";   C = "NSE:ANUP";        D = "NSE:HDFCAMC";  F = $null },
    @{ B = "NSE:COROMANDEL"; C = "NSE:AVANTIFEED";  D = "NSE:M&MFIN";   F = $null },
    @{ B = "NSE:EPIGRAL";    C = "NSE:COLPAL";      D = "NSE:MCX";      F = $null },
    @{ B = "NSE:HDFCSILVER"; C = "NSE:EMUDHRA";     D = $null;          F = $null },
    @{ B = "NSE:IPCALAB";    C = "NSE:HINDUNILVR";  D = $null;          F = $null },
    @{ B = "NSE:LEMONTREE";  C = "NSE:RADICO";      D = $null;          F = $null },
    @{ B = "NSE:MITTAL";     C = $null;             D = $null;          F = $null },
    @{ B = "NSE:RELCHEMQ";   C = $null;             D = $null;          F = $null }
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rec = $data[$i]

    if ($rec.B) { $ws.Cells.Item($row, 2).Value = $rec.B } else { $ws.Cells.Item($row, 2).Value = "" }
    if ($rec.C) { $ws.Cells.Item($row, 3).Value = $rec.C } else { $ws.Cells.Item($row, 3).Value = "" }
    if ($rec.D) { $ws.Cells.Item($row, 4).Value = $rec.D } else { $ws.Cells.Item($row, 4).Value = "" }
    $ws.Cells.Item($row, 5).Value = ""
    if ($rec.F) { $ws.Cells.Item($row, 6).Value = $rec.F } else { $ws.Cells.Item($row, 6).Value = "" }
}

# Remove the now-unused rows 11-14 entirely so the sheet dimension shrinks to A1:F10
$ws.Range("A11:F14").EntireRow.Delete() | Out-Null
